# Apply cell value updates to the cryptos list worksheet
# Generated from the authoritative cell-level diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.102.82"
$ws.Cells.Item(2, 5).Value = "  -1.41%  "
$ws.Cells.Item(3, 4).Value = "2.275.31"
$ws.Cells.Item(3, 5).Value = "  -0.43%  "
$ws.Cells.Item(4, 5).Value = "  -0.44%  "
$ws.Cells.Item(5, 4).Value = "'111.30"
$ws.Cells.Item(5, 5).Value = "  +1.18%  "
$ws.Cells.Item(6, 4).Value = "'264.31"
$ws.Cells.Item(6, 5).Value = "  -1.20%  "
$ws.Cells.Item(7, 4).Value = "'0.648"
$ws.Cells.Item(7, 5).Value = "  +3.90%  "
$ws.Cells.Item(8, 5).Value = "  -0.03%  "
$ws.Cells.Item(9, 4).Value = "'0.607"
$ws.Cells.Item(9, 5).Value = "  -1.04%  "
$ws.Cells.Item(10, 4).Value = "'46.45"
$ws.Cells.Item(10, 5).Value = "  -2.01%  "
$ws.Cells.Item(11, 4).Value = "'0.0934"
$ws.Cells.Item(11, 5).Value = "  -1.23%  "
$ws.Cells.Item(12, 4).Value = "'9.26"
$ws.Cells.Item(12, 5).Value = "  +3.10%  "
$ws.Cells.Item(13, 5).Value = "  +1.39%  "
$ws.Cells.Item(14, 4).Value = "'15.24"
$ws.Cells.Item(14, 5).Value = "  -2.80%  "
$ws.Cells.Item(15, 4).Value = "2.617.32"
$ws.Cells.Item(15, 5).Value = "  -0.45%  "
$ws.Cells.Item(16, 4).Value = "'0.856"
$ws.Cells.Item(16, 5).Value = "  +1.69%  "
$ws.Cells.Item(17, 4).Value = "2.285.17"
$ws.Cells.Item(17, 5).Value = "  -0.10%  "
$ws.Cells.Item(18, 4).Value = "43.170.01"
$ws.Cells.Item(18, 5).Value = "  -0.92%  "
$ws.Cells.Item(19, 5).Value = "  -0.90%  "
$ws.Cells.Item(20, 4).Value = "'6.72"
$ws.Cells.Item(20, 5).Value = "  -0.70%  "
$ws.Cells.Item(21, 4).Value = "'72.09"
$ws.Cells.Item(21, 5).Value = "  +0.01%  "
$ws.Cells.Item(22, 5).Value = "  -1.13%  "
$ws.Cells.Item(23, 4).Value = "'233.89"
$ws.Cells.Item(23, 5).Value = "  +0.79%  "
$ws.Cells.Item(24, 4).Value = "'2.87"
$ws.Cells.Item(24, 5).Value = "  +3.97%  "
$ws.Cells.Item(25, 4).Value = "'9.29"
$ws.Cells.Item(25, 5).Value = "  -4.34%  "
$ws.Cells.Item(26, 5).Value = "  +1.99%  "
$ws.Cells.Item(27, 4).Value = "'11.34"
$ws.Cells.Item(27, 5).Value = "  -2.57%  "
$ws.Cells.Item(28, 4).Value = "'40.82"
$ws.Cells.Item(28, 5).Value = "  -1.93%  "
$ws.Cells.Item(29, 5).Value = "  -1.07%  "
$ws.Cells.Item(30, 5).Value = "  -1.27%  "
$ws.Cells.Item(31, 4).Value = "'172.96"
$ws.Cells.Item(31, 5).Value = "  -1.58%  "
$ws.Cells.Item(32, 4).Value = "'21.40"
$ws.Cells.Item(32, 5).Value = "  -0.54%  "
$ws.Cells.Item(33, 4).Value = "'0.0895"
$ws.Cells.Item(33, 5).Value = "  -3.15%  "
$ws.Cells.Item(34, 4).Value = "'5.63"
$ws.Cells.Item(34, 5).Value = "  +0.37%  "
$ws.Cells.Item(35, 5).Value = "  +4.29%  "
$ws.Cells.Item(36, 4).Value = "'0.0378"
$ws.Cells.Item(36, 5).Value = "  +3.43%  "
$ws.Cells.Item(37, 4).Value = "'4.66"
$ws.Cells.Item(37, 5).Value = "  -1.04%  "
$ws.Cells.Item(38, 4).Value = "'3.90"
$ws.Cells.Item(38, 5).Value = "  +3.22%  "
$ws.Cells.Item(39, 5).Value = "  -2.87%  "
$ws.Cells.Item(40, 4).Value = "'2.57"
$ws.Cells.Item(40, 5).Value = "  +7.19%  "
$ws.Cells.Item(41, 5).Value = "  +3.43%  "
$ws.Cells.Item(42, 4).Value = "'74.27"
$ws.Cells.Item(42, 5).Value = "  +2.55%  "
$ws.Cells.Item(43, 5).Value = "  -2.65%  "
$ws.Cells.Item(44, 5).Value = "  -2.30%  "
$ws.Cells.Item(45, 5).Value = "  -0.11%  "
$ws.Cells.Item(46, 5).Value = "  -1.39%  "
$ws.Cells.Item(47, 2).Value = "TrustWalletToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(47, 4).Value = "'1.26"
$ws.Cells.Item(47, 5).Value = "  +3.69%  "
$ws.Cells.Item(48, 2).Value = "FraxShare"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(48, 4).Value = "'8.53"
$ws.Cells.Item(48, 5).Value = "  -3.31%  "
$ws.Cells.Item(49, 5).Value = "  +0.64%  "
$ws.Cells.Item(50, 4).Value = "'99.62"
$ws.Cells.Item(50, 5).Value = "  -2.37%  "
$ws.Cells.Item(51, 4).Value = "'0.601"
$ws.Cells.Item(51, 5).Value = "  +10.79%  "
